$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Table1 (RQ1: Cause of Flakiness?)  B3:C16  -- header row 3 unchanged
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "Tolerance"
$ws.Range("C4").Value = 10

$ws.Range("B5").Value = "Memory"
$ws.Range("C5").Value = 2

$ws.Range("B6").Value = "Dependencies / Environment"
$ws.Range("C6").Value = 2

$ws.Range("B7").Value = "Order of Events"
$ws.Range("C7").Value = 14

$ws.Range("B8").Value = "Concurrency"
$ws.Range("C8").Value = 3

$ws.Range("B9").Value = "Async Wait"
$ws.Range("C9").Value = 5

$ws.Range("B10").Value = "Delay"
$ws.Range("C10").Value = 2

# ---------------------------------------------------------------------------
# Table2 (RQ2: Fix for Flakiness?)  E3:F16  -- header row 3 unchanged
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = "(Tolerance) Increase acceptance levels in assert statements"
$ws.Range("F4").Value = 5

$ws.Range("E5").Value = "(Memory) deallocate unused space"
$ws.Range("F5").Value = 1

$ws.Range("E6").Value = "(Dependencies / Environment) update "
$ws.Range("F6").Value = 2

$ws.Range("E7").Value = "(Order of Events) implement event ordering"
$ws.Range("F7").Value = 3

$ws.Range("E8").Value = "(Tolerance) update conditionals"
$ws.Range("F8").Value = 3

$ws.Range("E9").Value = "(Order of Events) setup state"
$ws.Range("F9").Value = 2

$ws.Range("E10").Value = "(Order of Events) tear down state after shutdown"
$ws.Range("F10").Value = 9

$ws.Range("E11").Value = "(Tolerance) partial functionality"
$ws.Range("F11").Value = 2

$ws.Range("E12").Value = "(Concurrency) locks"
$ws.Range("F12").Value = 3

$ws.Range("E13").Value = "(Async Wait) added waitFor"
$ws.Range("F13").Value = 5

$ws.Range("E14").Value = "(Memory) added memory for test"
$ws.Range("F14").Value = 1

$ws.Range("E15").Value = "(Delay) add custom delay / wait"
$ws.Range("F15").Value = 2

# ---------------------------------------------------------------------------
# Table3 (Programming Language)  B18:C29  -- header row 18 unchanged
# ---------------------------------------------------------------------------
$ws.Range("B19").Value = "Scala"
$ws.Range("C19").Value = 1

$ws.Range("B20").Value = "Swift"
$ws.Range("C20").Value = 10

$ws.Range("B21").Value = "TypeScript"
$ws.Range("C21").Value = 3

$ws.Range("B22").Value = "Python"
$ws.Range("C22").Value = 4

$ws.Range("B23").Value = "JavaScript"
$ws.Range("C23").Value = 2

$ws.Range("B24").Value = "Java"
$ws.Range("C24").Value = 9

$ws.Range("B25").Value = "C++"
$ws.Range("C25").Value = 4

$ws.Range("B26").Value = "C#"
$ws.Range("C26").Value = 1

$ws.Range("B27").Value = "Go"
$ws.Range("C27").Value = 3

$ws.Range("B28").Value = "Kotlin"
$ws.Range("C28").Value = 1

# ---------------------------------------------------------------------------
# Sheet view: scroll position + selection
# ---------------------------------------------------------------------------
$ws.Range("A9").Select()
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("C27").Select()
